# The post "「ひとを幸せにするには？」" (row 342) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one and
# automatically shrinks the used range (dimension) from A1:C487 to A1:C486.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(342).Delete()
